# Apply the "output folder feature" edits described in the commit:
#  - Input sheet: record the source data filename in a new context column (D)
#  - Calc sheet: add two new "Erfolgsrate" (success-rate) columns (BG/BH)
#  - Slightly revised Monte-Carlo derived values that ripple into Calc/Results
#  - Constants sheet: minor precision tweak on one constant

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Input
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

# Store the full path of the source .dat file alongside the sample record.
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"

# Widen column D so the filename is fully visible.
$wsInput.Columns.Item(4).ColumnWidth = 65.8333333333

# ---------------------------------------------------------------------------
# Sheet: Calc
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calc")

# Updated Monte-Carlo derived figures (re-run with the new output-folder logic).
$wsCalc.Range("AP3").Value = 0.5231
$wsCalc.Range("AQ3").Value = 0.185082708724893
$wsCalc.Range("AW3").Value = 0.5401
$wsCalc.Range("AX3").Value = 0.5354463754901183
$wsCalc.Range("AY3").Value = 0.1911333382405399
$wsCalc.Range("BC3").Value = 0.5693696230851605
$wsCalc.Range("BE3").Value = 267.7231877450591
$wsCalc.Range("BF3").Value = 0.1894864899022846

# New columns: Monte-Carlo error "success rate" (uncorrected / corrected).
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BG3").Value = 100
$wsCalc.Range("BH3").Value = 100

# Column width adjustments.
$wsCalc.Columns.Item(49).ColumnWidth = 8.833333333333332
$wsCalc.Columns.Item(51).ColumnWidth = 19.833333333333336
$wsCalc.Columns.Item(57).ColumnWidth = 18.833333333333336
$wsCalc.Columns.Item(59).ColumnWidth = 31.833333333333332
$wsCalc.Columns.Item(60).ColumnWidth = 29.833333333333336

# ---------------------------------------------------------------------------
# Sheet: Results
# ---------------------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")

$wsResults.Range("N3").Value = 0.5231
$wsResults.Range("P3").Value = 0.5401
$wsResults.Range("R3").Value = 0.5693696230851605

$wsResults.Columns.Item(16).ColumnWidth = 7.833333333333333

# ---------------------------------------------------------------------------
# Sheet: Constants
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")

$wsConstants.Range("B3").Value = 0.00005
